$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1978021978021978
$ws.Range("C2").Value = 0.5311355311355311
$ws.Range("J2").Value = 0.01098901098901099
$ws.Range("P2").Value = 0.1391941391941392
$ws.Range("S2").Value = 0.1208791208791209
$ws.Range("B3").Value = 0.01910828025477707
$ws.Range("C3").Value = 0.03184713375796178
$ws.Range("J3").Value = 0.03184713375796178
$ws.Range("P3").Value = 0.6878980891719745
$ws.Range("S3").Value = 0.2292993630573248
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.07766990291262135
$ws.Range("D6").Value = 0.004854368932038835
$ws.Range("F6").Value = 0.04854368932038835
$ws.Range("J6").Value = 0.2233009708737864
$ws.Range("O6").Value = 0.03398058252427184
$ws.Range("Q6").Value = 0.1650485436893204
$ws.Range("R6").Value = 0.1359223300970874
$ws.Range("S6").Value = 0.3106796116504854
$ws.Range("B7").Value = 0.08552631578947369
$ws.Range("D7").Value = 0.006578947368421052
$ws.Range("F7").Value = 0.08552631578947369
$ws.Range("J7").Value = 0.1513157894736842
$ws.Range("O7").Value = 0.0131578947368421
$ws.Range("Q7").Value = 0.131578947368421
$ws.Range("R7").Value = 0.1118421052631579
$ws.Range("S7").Value = 0.4144736842105263
$ws.Range("B8").Value = 0.09487179487179487
$ws.Range("D8").Value = 0.01794871794871795
$ws.Range("F8").Value = 0.07179487179487179
$ws.Range("J8").Value = 0.1128205128205128
$ws.Range("O8").Value = 0.02307692307692308
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.1256410256410256
$ws.Range("S8").Value = 0.3871794871794872
$ws.Range("B9").Value = 0.07926829268292683
$ws.Range("D9").Value = 0.01219512195121951
$ws.Range("F9").Value = 0.05487804878048781
$ws.Range("J9").Value = 0.1768292682926829
$ws.Range("O9").Value = 0.006097560975609756
$ws.Range("Q9").Value = 0.1707317073170732
$ws.Range("R9").Value = 0.0975609756097561
$ws.Range("S9").Value = 0.4024390243902439
$ws.Range("B10").Value = 0.1130653266331658
$ws.Range("D10").Value = 0.01842546063651591
$ws.Range("E10").Value = 0.002512562814070352
$ws.Range("F10").Value = 0.076214405360134
$ws.Range("J10").Value = 0.1432160804020101
$ws.Range("O10").Value = 0.01675041876046901
$ws.Range("Q10").Value = 0.1876046901172529
$ws.Range("R10").Value = 0.1038525963149079
$ws.Range("S10").Value = 0.338358458961474
$ws.Range("G11").Value = 0.09954751131221719
$ws.Range("J11").Value = 0.08144796380090498
$ws.Range("K11").Value = 0.16289592760181
$ws.Range("L11").Value = 0.6561085972850679
$ws.Range("G12").Value = 0.7651006711409396
$ws.Range("J12").Value = 0.174496644295302
$ws.Range("K12").Value = 0.006711409395973154
$ws.Range("L12").Value = 0.04697986577181208
$ws.Range("S12").Value = 0.006711409395973154
$ws.Range("G13").Value = 0.7241379310344828
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("F15").Value = 0.02564102564102564
$ws.Range("H15").Value = 0.1794871794871795
$ws.Range("I15").Value = 0.05641025641025641
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.01025641025641026
$ws.Range("O15").Value = 0.04102564102564103
$ws.Range("S15").Value = 0.2205128205128205
$ws.Range("F16").Value = 0.01807228915662651
$ws.Range("H16").Value = 0.2650602409638554
$ws.Range("I16").Value = 0.1144578313253012
$ws.Range("J16").Value = 0.3433734939759036
$ws.Range("K16").Value = 0.09036144578313253
$ws.Range("M16").Value = 0.01807228915662651
$ws.Range("O16").Value = 0.04819277108433735
$ws.Range("S16").Value = 0.1024096385542169
$ws.Range("F17").Value = 0.008086253369272238
$ws.Range("H17").Value = 0.1563342318059299
$ws.Range("I17").Value = 0.1078167115902965
$ws.Range("J17").Value = 0.4582210242587601
$ws.Range("K17").Value = 0.07277628032345014
$ws.Range("M17").Value = 0.01347708894878706
$ws.Range("O17").Value = 0.07547169811320754
$ws.Range("S17").Value = 0.1078167115902965
$ws.Range("F18").Value = 0.02553191489361702
$ws.Range("H18").Value = 0.1574468085106383
$ws.Range("I18").Value = 0.09787234042553192
$ws.Range("J18").Value = 0.425531914893617
$ws.Range("K18").Value = 0.09787234042553192
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("N18").Value = 0.00425531914893617
$ws.Range("O18").Value = 0.06808510638297872
$ws.Range("S18").Value = 0.1021276595744681
$ws.Range("F19").Value = 0.01253616200578592
$ws.Range("H19").Value = 0.2121504339440694
$ws.Range("I19").Value = 0.0703953712632594
$ws.Range("J19").Value = 0.4146576663452266
$ws.Range("K19").Value = 0.1002892960462874
$ws.Range("M19").Value = 0.01735776277724204
$ws.Range("O19").Value = 0.07521697203471553
$ws.Range("S19").Value = 0.09739633558341369
